$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SQL server IP value in E2 (SqlIP column)
$ws.Range("E2").Value = "192.168.0.24"

# Move the active selection to H6, matching the recorded UI state
$ws.Range("H6").Select()
